# Applies the "add HideFile field" change to TestDataSet / TechSet tables,
# fills in the new HideFile data, and restores the active-sheet/selection state.

$wb = $excel.ActiveWorkbook

$wsTestDataSet = $wb.Worksheets.Item(1)   # TestDataSet
$wsTechSet     = $wb.Worksheets.Item(2)   # TechSet
$wsFileSet     = $wb.Worksheets.Item(3)   # FileSet

# Remember a representative column width so the new columns match the others.
$stdWidth = $wsTestDataSet.Columns.Item(1).ColumnWidth

## --- TestDataSet: add the HideFile column to the table -------------------
$loTestDataSet = $wsTestDataSet.ListObjects.Item(1)
$loTestDataSet.ListColumns.Add() | Out-Null
$wsTestDataSet.Range("K1").Value = "HideFile"
$wsTestDataSet.Columns.Item(11).ColumnWidth = $stdWidth

# Fill in the HideFile values for the data rows (row 2 stays blank).
$wsTestDataSet.Range("K3").Value = 1
$wsTestDataSet.Range("K4").Value = 1
$wsTestDataSet.Range("K5").Value = 1
$wsTestDataSet.Range("K6").Value = 1
$wsTestDataSet.Range("K7").Value = 1
$wsTestDataSet.Range("K8").Value = 1
$wsTestDataSet.Range("K9").Value = 1
$wsTestDataSet.Range("K10").Value = 1
$wsTestDataSet.Range("K11").Value = 1

## --- TechSet: add the HideFile column to the table ------------------------
$loTechSet = $wsTechSet.ListObjects.Item(1)
$loTechSet.ListColumns.Add() | Out-Null
$wsTechSet.Range("E1").Value = "HideFile"
$wsTechSet.Columns.Item(5).ColumnWidth = $stdWidth

# Fill in HideFile values: Status 00 -> False, 01 -> True, 02 -> True.
$wsTechSet.Range("E2").Value = $false
$wsTechSet.Range("E3").Value = $true
$wsTechSet.Range("E4").Value = $true

## --- Restore view/selection state -----------------------------------------
# FileSet previously held the active selection; clear it back to A1.
$wsFileSet.Activate()
$wsFileSet.Range("A1").Select()

# TestDataSet becomes the active sheet/tab again, with K10 selected.
$wsTestDataSet.Activate()
$wsTestDataSet.Range("K10").Select()
